{"js": "// Insert three new bullet paragraphs into the \"PARTNER - Siege Analytics\"\n// section, right after the \"Research & Data Analytics Leadership\" heading\n// paragraph and before the existing \"\u2022 Conceived, architected, ...\" bullet.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the anchor paragraph: \"Research & Data Analytics Leadership\"\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Research & Data Analytics Leadership\") {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find the 'Research & Data Analytics Leadership' paragraph\");\n}\n\n// Insert in reverse order, each time right \"After\" the anchor paragraph,\n// so the final reading order is:\n//   Research & Data Analytics Leadership\n//   \u2022 Uncovered decades of demographic miscoding ...\n//   \u2022 Developed Python boundary estimation algorithm ...\n//   \u2022 Algorithm reduced mapping costs by 75% ...\n//   \u2022 Conceived, architected, engineered and deployed ... (pre-existing)\nanchor.insertParagraph(\n  \"\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations $5M+ and enabling smaller nonprofits to conduct redistricting analysis\",\n  \"After\"\n);\nanchor.insertParagraph(\n  \"\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n  \"After\"\n);\nanchor.insertParagraph(\n  \"\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# Insert three new bullet paragraphs into the \"PARTNER - Siege Analytics\"\n# section, right after the \"Research & Data Analytics Leadership\" heading\n# paragraph and before the existing \"- Conceived, architected, ...\" bullet.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (\"Research & Data Analytics Leadership\") robustly via Find.\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.Text = 'Research & Data Analytics Leadership'\n$found = $find.Execute()\nif (-not $found) {\n    throw 'Anchor paragraph \"Research & Data Analytics Leadership\" not found'\n}\n\n# Resolve which paragraph (by 1-based index) the found range falls within.\n$anchorStart = $searchRange.Start\n$paraCount = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $paraCount; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ($candidate.Range.Start -le $anchorStart -and $anchorStart -lt $candidate.Range.End) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw 'Could not resolve the anchor paragraph index'\n}\n\n# New bullet paragraphs to insert, in the order they should appear.\n$newBullets = @(\n    '\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters',\n    '\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States',\n    '\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations $5M+ and enabling smaller nonprofits to conduct redistricting analysis'\n)\n\nforeach ($bulletText in $newBullets) {\n    $anchorRange = $d.Paragraphs.Item($anchorIndex).Range\n    $anchorRange.Collapse(0)   # wdCollapseEnd\n    $anchorRange.InsertParagraphAfter()\n    $newParaRange = $d.Paragraphs.Item($anchorIndex + 1).Range\n    $newParaRange.Text = $bulletText\n    $anchorIndex = $anchorIndex + 1\n}\n"}
